# Split Course Code and Level on Summer Reporting File Spec
#
# Column G used to hold "Ministry Course Code and Level" (e.g. "ENST 12").
# It is being split into two columns:
#   G: "Ministry Course Code"  (e.g. "ENST")
#   H: "Ministry Course Level" (e.g. 12, numeric)
# Every column that was to the right of the old G (Session Date, Final
# Percent, Final Letter Grade, Credits) shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh column at H - this pushes the old H:K (Session Date,
# Final Percent, Final Letter Grade, Credits) to I:L and gives us an
# empty column to hold the new "Ministry Course Level" values.
$ws.Columns("H").Insert()

# Re-label column G and fill in the course-code-only values.
$ws.Range("G1").Value = "Ministry Course Code"

$ws.Range("G2").Value = "ENST"
$ws.Range("G3").Value = "ENST"
$ws.Range("G4").Value = "ENST"

# Label the new column and fill in the course-level values.
$ws.Range("H1").Value = "Ministry Course Level"

$ws.Range("H2").Value = 12
$ws.Range("H3").Value = 12
$ws.Range("H4").Value = 12

# Leave the new split columns selected, matching the author's selection
# after performing the split.
$ws.Columns("G:H").Select()
